$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the question text: B1 and C1 both hold the same question now
$ws.Range("B1").Value = "Which is the capital city of Karnataka ? "
$ws.Range("C1").Value = "Which is the capital city of Karnataka ? "

# Remove the old answer in B2, and set C2 to a single blank space
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = " "

# Widen columns B and C to fit the longer question text
$ws.Columns("B:C").ColumnWidth = 40.166666666666664
